# DDAS_Upload_CaseStudy_5.xlsx
#
# QCFailed and QCPassed rows have been replaced with QCCompleted rows: two
# brand-new QC input rows (Huda, Syed and Liang, Cheng Yi) are added as the
# new top data rows, pushing the pre-existing data down. The sheet view is
# also scrolled/selected to K2 so the new rows are visible.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two new blank rows above the current row 2 (the first data row).
#    This shifts: old row2 (Jones, Sheena) -> row4
#                 old row3 (blank)         -> row5
#                 old row4 (Mok)           -> row6
#                 old row5 (Jolly)         -> row7
#                 old row6 (Xing)          -> row8
#                 old row7 (Van Der Heijde)-> row9
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

# ---------------------------------------------------------------------------
# 2. The two newly inserted rows (2:3) come back with a generic style; copy
#    the per-column formatting from row 4 (the shifted-down, still blank
#    placeholder row) onto them so every column keeps its original look.
# ---------------------------------------------------------------------------
$ws.Range("A4:S4").Copy()
$ws.Range("A2:S3").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. The data that used to sit in row 2 (Jones, Sheena) needs to end up
#    below the blank placeholder row again, i.e. in row 5, with row 4 left
#    empty - matching how the sheet looked before (data row, then blank
#    separator row) just shifted down by two rows.
# ---------------------------------------------------------------------------
$ws.Range("A4:S4").Copy()
$ws.Range("A5").PasteSpecial(-4163)      # xlPasteValues
$excel.CutCopyMode = 0
$ws.Range("A4:S4").ClearContents()

# ---------------------------------------------------------------------------
# 4. Populate the new row 2 with the Huda, Syed QC record.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "PI"
$ws.Range("B2").Value = "0000/0004"
$ws.Range("C2").Value = "0000/0000"
$ws.Range("D2").Value = 9951
$ws.Range("F2").Value = "Huda, Syed"
$ws.Range("I2").Value = "Syed"
$ws.Range("K2").Value = "Huda"
$ws.Range("L2").Value = "St. Petersburg"

# ---------------------------------------------------------------------------
# 5. Populate the new row 3 with the Liang, Cheng Yi QC record.
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "PI"
$ws.Range("B3").Value = "0000/0004"
$ws.Range("C3").Value = "0000/0000"
$ws.Range("D3").Value = 2274
$ws.Range("F3").Value = "Liang, Cheng Yi"
$ws.Range("I3").Value = "Yi"
$ws.Range("J3").Value = "Liang"
$ws.Range("K3").Value = "Cheng"

# ---------------------------------------------------------------------------
# 6. Update the sheet view to scroll/select K2, as in the saved workbook.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollColumn = 11
$ws.Range("K2").Select()
